# Rename transcript speaker tag "R1" to "T" in column D for the specified rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(19, 21, 24, 25, 26, 28, 31, 32, 35, 41, 52, 55)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Value = "T"
}
